$d = $word.ActiveDocument

$d.Content.Find.Execute("477÷5=95, 2", $true, $false, $false, $false, $false, $true, 1, $false, "737÷6=122, 5", 2) | Out-Null
$d.Content.Find.Execute("869÷6=144, 5", $true, $false, $false, $false, $false, $true, 1, $false, "938÷8=117, 2", 2) | Out-Null
$d.Content.Find.Execute("837÷2=418, 1", $true, $false, $false, $false, $false, $true, 1, $false, "727÷2=363, 1", 2) | Out-Null
$d.Content.Find.Execute("935÷2=467, 1", $true, $false, $false, $false, $false, $true, 1, $false, "880÷8=110, 0", 2) | Out-Null
$d.Content.Find.Execute("535÷4=133, 3", $true, $false, $false, $false, $false, $true, 1, $false, "108÷8=13, 4", 2) | Out-Null
$d.Content.Find.Execute("702÷8=87, 6", $true, $false, $false, $false, $false, $true, 1, $false, "892÷8=111, 4", 2) | Out-Null
$d.Content.Find.Execute("505÷5=101, 0", $true, $false, $false, $false, $false, $true, 1, $false, "669÷3=223, 0", 2) | Out-Null
$d.Content.Find.Execute("777÷6=129, 3", $true, $false, $false, $false, $false, $true, 1, $false, "940÷5=188, 0", 2) | Out-Null
$d.Content.Find.Execute("420÷9=46, 6", $true, $false, $false, $false, $false, $true, 1, $false, "980÷8=122, 4", 2) | Out-Null
$d.Content.Find.Execute("183÷4=45, 3", $true, $false, $false, $false, $false, $true, 1, $false, "327÷4=81, 3", 2) | Out-Null
$d.Content.Find.Execute("661÷6=110, 1", $true, $false, $false, $false, $false, $true, 1, $false, "891÷2=445, 1", 2) | Out-Null
$d.Content.Find.Execute("453÷6=75, 3", $true, $false, $false, $false, $false, $true, 1, $false, "283÷4=70, 3", 2) | Out-Null
$d.Content.Find.Execute("635÷9=70, 5", $true, $false, $false, $false, $false, $true, 1, $false, "677÷5=135, 2", 2) | Out-Null
$d.Content.Find.Execute("672÷7=96, 0", $true, $false, $false, $false, $false, $true, 1, $false, "663÷9=73, 6", 2) | Out-Null
$d.Content.Find.Execute("489÷8=61, 1", $true, $false, $false, $false, $false, $true, 1, $false, "396÷7=56, 4", 2) | Out-Null
$d.Content.Find.Execute("730÷5=146, 0", $true, $false, $false, $false, $false, $true, 1, $false, "813÷5=162, 3", 2) | Out-Null
$d.Content.Find.Execute("916÷9=101, 7", $true, $false, $false, $false, $false, $true, 1, $false, "681÷5=136, 1", 2) | Out-Null
$d.Content.Find.Execute("868÷8=108, 4", $true, $false, $false, $false, $false, $true, 1, $false, "465÷7=66, 3", 2) | Out-Null
$d.Content.Find.Execute("287÷6=47, 5", $true, $false, $false, $false, $false, $true, 1, $false, "990÷6=165, 0", 2) | Out-Null
$d.Content.Find.Execute("152÷3=50, 2", $true, $false, $false, $false, $false, $true, 1, $false, "153÷8=19, 1", 2) | Out-Null
$d.Content.Find.Execute("693÷9=77, 0", $true, $false, $false, $false, $false, $true, 1, $false, "771÷6=128, 3", 2) | Out-Null
$d.Content.Find.Execute("514÷8=64, 2", $true, $false, $false, $false, $false, $true, 1, $false, "337÷7=48, 1", 2) | Out-Null
$d.Content.Find.Execute("373÷7=53, 2", $true, $false, $false, $false, $false, $true, 1, $false, "314÷4=78, 2", 2) | Out-Null
$d.Content.Find.Execute("519÷2=259, 1", $true, $false, $false, $false, $false, $true, 1, $false, "633÷3=211, 0", 2) | Out-Null
$d.Content.Find.Execute("811÷4=202, 3", $true, $false, $false, $false, $false, $true, 1, $false, "346÷9=38, 4", 2) | Out-Null
